$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": status columns (zh-cn / de-de) flip from "Ready for
# handoff" to "Handed back: in sync with en-US" now that the report has run.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the status columns so the longer text fits.
$wsOverview.Range("E1").ColumnWidth = 29.17
$wsOverview.Range("F1").ColumnWidth = 29.17

# ---------------------------------------------------------------------------
# Sheet "zh-cn": fill in the handback info (target file / handback file /
# handback datetime) for both rows now that a handback report was generated.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/167bf863af857649f39fbff57d88507ff513d512/e2e/3b97f684-42b1-4dc9-b1f3-c2d704a2ccbe.md", "", "", "3b97f684-42b1-4dc9-b1f3-c2d704a2ccbe.md")
$wsZhCn.Range("J2").Value = "3b97f684-42b1-4dc9-b1f3-c2d704a2ccbe.c7d4381cf837edcb9b78371d86387f330e685a24.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-31 10:33:19"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/167bf863af857649f39fbff57d88507ff513d512/e2e/3daaecf5-41be-4124-a918-8a941a174150.md", "", "", "3daaecf5-41be-4124-a918-8a941a174150.md")
$wsZhCn.Range("J3").Value = "3daaecf5-41be-4124-a918-8a941a174150.3c9df6919bf070db1a62de59649d2c317b19bd6d.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-31 10:33:19"

$wsZhCn.Range("C1").ColumnWidth = 29.17
$wsZhCn.Range("I1").ColumnWidth = 39.17
$wsZhCn.Range("J1").ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# Sheet "de-de": same handback info, but the de-de xlf names + its own
# handback datetime.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/167bf863af857649f39fbff57d88507ff513d512/e2e/3b97f684-42b1-4dc9-b1f3-c2d704a2ccbe.md", "", "", "3b97f684-42b1-4dc9-b1f3-c2d704a2ccbe.md")
$wsDeDe.Range("J2").Value = "3b97f684-42b1-4dc9-b1f3-c2d704a2ccbe.c7d4381cf837edcb9b78371d86387f330e685a24.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-31 10:33:26"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/167bf863af857649f39fbff57d88507ff513d512/e2e/3daaecf5-41be-4124-a918-8a941a174150.md", "", "", "3daaecf5-41be-4124-a918-8a941a174150.md")
$wsDeDe.Range("J3").Value = "3daaecf5-41be-4124-a918-8a941a174150.3c9df6919bf070db1a62de59649d2c317b19bd6d.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-31 10:33:26"

$wsDeDe.Range("C1").ColumnWidth = 29.17
$wsDeDe.Range("I1").ColumnWidth = 39.17
$wsDeDe.Range("J1").ColumnWidth = 39.17
